$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (pushes existing data rows down by one,
# carrying their formatting with them) and drop the row that now falls
# off the bottom of the table (old row 20 -> now row 21) so the sheet
# stays a fixed 20-row table (A1:Q20).
$ws.Rows("2:2").Insert()
$ws.Rows("21:21").Delete()

# Populate the new top data row with the latest IPO record.
$ws.Range("A2").Value = "2024-06-17"
$ws.Range("B2").Value = "라메디텍"
$ws.Range("C2").Value = "코스닥"
$ws.Range("D2").Value = 207.68
$ws.Range("E2").Value = "대신"
$ws.Range("F2").Value = 207.68
$ws.Range("G2").Value = "-"
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = "-"
$ws.Range("J2").Value = "-"
$ws.Range("K2").Value = "대표"
$ws.Range("L2").Value = "-"
$ws.Range("M2").Value = 16000
$ws.Range("N2").Value = 100
$ws.Range("O2").Value = "2024-06-05"
$ws.Range("P2").Value = "2024-06-11"
$ws.Range("Q2").Value = 957220
